# Separate game entity and model
# Add a new row (62) to the action/commit log sheet, mirroring the
# formatting of the existing "section header" rows (bold, column A + B).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = "Separate mongodb entities and pure chess models"
$ws.Range("A62").Font.Bold = $true
$ws.Range("B62").Value = "x"

# Update the active selection to match the new last cell, like Excel
# would after entering data in that row.
$excel.Goto($ws.Range("C62"))
